$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns before column D to make room for the two newest
# fiscal quarters (the sheet now reports 10 quarters instead of 8).
$ws.Columns("D:E").Insert()

# The newly inserted D:E columns have no formatting yet; clone it from the
# (now shifted) F:M columns so every row keeps its original date / number
# style instead of the generic column default.
$ws.Range("F5:M102").Copy()
$ws.Range("D5:E102").PasteSpecial(-4122)

$ws.Cells.Item(7, 4).Value = 43465
$ws.Cells.Item(7, 5).Value = 43373
$ws.Cells.Item(7, 6).Value = 43281
$ws.Cells.Item(7, 7).Value = 43190
$ws.Cells.Item(7, 8).Value = 43100
$ws.Cells.Item(7, 9).Value = 43008
$ws.Cells.Item(7, 10).Value = 42916
$ws.Cells.Item(7, 11).Value = 42825
$ws.Cells.Item(7, 12).Value = 42735
$ws.Cells.Item(7, 13).Value = 42643
$ws.Cells.Item(8, 4).Value = 2064000
$ws.Cells.Item(8, 5).Value = 2212000
$ws.Cells.Item(8, 6).Value = 2156000
$ws.Cells.Item(8, 7).Value = 1832000
$ws.Cells.Item(8, 8).Value = 1942000
$ws.Cells.Item(8, 9).Value = 2110000
$ws.Cells.Item(8, 10).Value = 2102000
$ws.Cells.Item(8, 11).Value = 1740000
$ws.Cells.Item(8, 12).Value = 1524000
$ws.Cells.Item(8, 13).Value = 1566000
$ws.Cells.Item(9, 4).Value = 1003000
$ws.Cells.Item(9, 5).Value = 993000
$ws.Cells.Item(9, 6).Value = 960000
$ws.Cells.Item(9, 7).Value = 897000
$ws.Cells.Item(9, 8).Value = 1277000
$ws.Cells.Item(9, 9).Value = 1179000
$ws.Cells.Item(9, 10).Value = 1104000
$ws.Cells.Item(9, 11).Value = 1209000
$ws.Cells.Item(9, 12).Value = 927000
$ws.Cells.Item(9, 13).Value = 837000
$ws.Cells.Item(10, 4).Value = 1061000
$ws.Cells.Item(10, 5).Value = 1219000
$ws.Cells.Item(10, 6).Value = 1196000
$ws.Cells.Item(10, 7).Value = 935000
$ws.Cells.Item(10, 8).Value = 665000
$ws.Cells.Item(10, 9).Value = 931000
$ws.Cells.Item(10, 10).Value = 998000
$ws.Cells.Item(10, 11).Value = 531000
$ws.Cells.Item(10, 12).Value = 597000
$ws.Cells.Item(10, 13).Value = 729000
$ws.Cells.Item(12, 4).Value = "NA"
$ws.Cells.Item(12, 5).Value = "NA"
$ws.Cells.Item(12, 6).Value = "NA"
$ws.Cells.Item(12, 7).Value = "NA"
$ws.Cells.Item(12, 8).Value = "NA"
$ws.Cells.Item(12, 9).Value = "NA"
$ws.Cells.Item(12, 10).Value = "NA"
$ws.Cells.Item(12, 11).Value = "NA"
$ws.Cells.Item(12, 12).Value = "NA"
$ws.Cells.Item(12, 13).Value = "NA"
$ws.Cells.Item(13, 4).Value = 0
$ws.Cells.Item(13, 5).Value = 0
$ws.Cells.Item(13, 6).Value = 0
$ws.Cells.Item(13, 7).Value = 0
$ws.Cells.Item(13, 8).Value = 0
$ws.Cells.Item(13, 9).Value = 0
$ws.Cells.Item(13, 10).Value = 0
$ws.Cells.Item(13, 11).Value = 0
$ws.Cells.Item(13, 12).Value = 0
$ws.Cells.Item(13, 13).Value = 0
$ws.Cells.Item(14, 4).Value = 20000
$ws.Cells.Item(14, 5).Value = 22000
$ws.Cells.Item(14, 6).Value = 39000
$ws.Cells.Item(14, 7).Value = 31000
$ws.Cells.Item(14, 8).Value = 30000
$ws.Cells.Item(14, 9).Value = 23000
$ws.Cells.Item(14, 10).Value = 24000
$ws.Cells.Item(14, 11).Value = 39000
$ws.Cells.Item(14, 12).Value = 81000
$ws.Cells.Item(14, 13).Value = 22000
$ws.Cells.Item(15, 4).Value = 108000
$ws.Cells.Item(15, 5).Value = 99000
$ws.Cells.Item(15, 6).Value = 97000
$ws.Cells.Item(15, 7).Value = 94000
$ws.Cells.Item(15, 8).Value = 97000
$ws.Cells.Item(15, 9).Value = 95000
$ws.Cells.Item(15, 10).Value = 90000
$ws.Cells.Item(15, 11).Value = 90000
$ws.Cells.Item(15, 12).Value = 82000
$ws.Cells.Item(15, 13).Value = 101000
$ws.Cells.Item(17, 4).Value = 2018000
$ws.Cells.Item(17, 5).Value = 1915000
$ws.Cells.Item(17, 6).Value = 1885000
$ws.Cells.Item(17, 7).Value = 1803000
$ws.Cells.Item(17, 8).Value = 1805000
$ws.Cells.Item(17, 9).Value = 1683000
$ws.Cells.Item(17, 10).Value = 1615000
$ws.Cells.Item(17, 11).Value = 1583000
$ws.Cells.Item(17, 12).Value = 1283000
$ws.Cells.Item(17, 13).Value = 1166000
$ws.Cells.Item(18, 4).Value = 46000
$ws.Cells.Item(18, 5).Value = 297000
$ws.Cells.Item(18, 6).Value = 271000
$ws.Cells.Item(18, 7).Value = 29000
$ws.Cells.Item(18, 8).Value = 137000
$ws.Cells.Item(18, 9).Value = 427000
$ws.Cells.Item(18, 10).Value = 487000
$ws.Cells.Item(18, 11).Value = 157000
$ws.Cells.Item(18, 12).Value = 241000
$ws.Cells.Item(18, 13).Value = 400000
$ws.Cells.Item(20, 4).Value = 6000
$ws.Cells.Item(20, 5).Value = 4000
$ws.Cells.Item(20, 6).Value = 9000
$ws.Cells.Item(20, 7).Value = -4000
$ws.Cells.Item(20, 8).Value = 11000
$ws.Cells.Item(20, 9).Value = 11000
$ws.Cells.Item(20, 10).Value = 9000
$ws.Cells.Item(20, 11).Value = 6000
$ws.Cells.Item(20, 12).Value = 8000
$ws.Cells.Item(20, 13).Value = 7000
$ws.Cells.Item(21, 4).Value = 160000
$ws.Cells.Item(21, 5).Value = 400000
$ws.Cells.Item(21, 6).Value = 377000
$ws.Cells.Item(21, 7).Value = 119000
$ws.Cells.Item(21, 8).Value = 245000
$ws.Cells.Item(21, 9).Value = 533000
$ws.Cells.Item(21, 10).Value = 586000
$ws.Cells.Item(21, 11).Value = 253000
$ws.Cells.Item(21, 12).Value = 331000
$ws.Cells.Item(21, 13).Value = 508000
$ws.Cells.Item(22, 4).Value = 16000
$ws.Cells.Item(22, 5).Value = 17000
$ws.Cells.Item(22, 6).Value = 21000
$ws.Cells.Item(22, 7).Value = 19000
$ws.Cells.Item(22, 8).Value = 22000
$ws.Cells.Item(22, 9).Value = 21000
$ws.Cells.Item(22, 10).Value = 22000
$ws.Cells.Item(22, 11).Value = 21000
$ws.Cells.Item(22, 12).Value = 18000
$ws.Cells.Item(22, 13).Value = 5000
$ws.Cells.Item(23, 4).Value = 36000
$ws.Cells.Item(23, 5).Value = 284000
$ws.Cells.Item(23, 6).Value = 259000
$ws.Cells.Item(23, 7).Value = 6000
$ws.Cells.Item(23, 8).Value = 126000
$ws.Cells.Item(23, 9).Value = 417000
$ws.Cells.Item(23, 10).Value = 474000
$ws.Cells.Item(23, 11).Value = 142000
$ws.Cells.Item(23, 12).Value = 231000
$ws.Cells.Item(23, 13).Value = 402000
$ws.Cells.Item(24, 4).Value = 13000
$ws.Cells.Item(24, 5).Value = 67000
$ws.Cells.Item(24, 6).Value = 66000
$ws.Cells.Item(24, 7).Value = 2000
$ws.Cells.Item(24, 8).Value = 48000
$ws.Cells.Item(24, 9).Value = 158000
$ws.Cells.Item(24, 10).Value = 181000
$ws.Cells.Item(24, 11).Value = 49000
$ws.Cells.Item(24, 12).Value = 117000
$ws.Cells.Item(24, 13).Value = 146000
$ws.Cells.Item(25, 4).Value = 0
$ws.Cells.Item(25, 5).Value = 0
$ws.Cells.Item(25, 6).Value = 0
$ws.Cells.Item(25, 7).Value = 0
$ws.Cells.Item(25, 8).Value = 0
$ws.Cells.Item(25, 9).Value = 0
$ws.Cells.Item(25, 10).Value = 0
$ws.Cells.Item(25, 11).Value = 0
$ws.Cells.Item(25, 12).Value = 0
$ws.Cells.Item(25, 13).Value = 0
$ws.Cells.Item(26, 4).Value = 23000
$ws.Cells.Item(26, 5).Value = 217000
$ws.Cells.Item(26, 6).Value = 193000
$ws.Cells.Item(26, 7).Value = 4000
$ws.Cells.Item(26, 8).Value = 78000
$ws.Cells.Item(26, 9).Value = 259000
$ws.Cells.Item(26, 10).Value = 293000
$ws.Cells.Item(26, 11).Value = 93000
$ws.Cells.Item(26, 12).Value = 114000
$ws.Cells.Item(26, 13).Value = 256000
$ws.Cells.Item(27, 4).Value = 23000
$ws.Cells.Item(27, 5).Value = 217000
$ws.Cells.Item(27, 6).Value = 193000
$ws.Cells.Item(27, 7).Value = 4000
$ws.Cells.Item(27, 8).Value = 78000
$ws.Cells.Item(27, 9).Value = 259000
$ws.Cells.Item(27, 10).Value = 293000
$ws.Cells.Item(27, 11).Value = 93000
$ws.Cells.Item(27, 12).Value = 114000
$ws.Cells.Item(27, 13).Value = 256000
$ws.Cells.Item(28, 4).Value = 0
$ws.Cells.Item(28, 5).Value = 0
$ws.Cells.Item(28, 6).Value = 0
$ws.Cells.Item(28, 7).Value = 0
$ws.Cells.Item(28, 8).Value = 0
$ws.Cells.Item(28, 9).Value = 0
$ws.Cells.Item(28, 10).Value = 0
$ws.Cells.Item(28, 11).Value = 0
$ws.Cells.Item(28, 12).Value = 0
$ws.Cells.Item(28, 13).Value = 0
$ws.Cells.Item(29, 4).Value = 0
$ws.Cells.Item(29, 5).Value = "NA"
$ws.Cells.Item(29, 6).Value = "NA"
$ws.Cells.Item(29, 7).Value = "NA"
$ws.Cells.Item(29, 8).Value = 237000
$ws.Cells.Item(29, 9).Value = "NA"
$ws.Cells.Item(29, 10).Value = "NA"
$ws.Cells.Item(29, 11).Value = "NA"
$ws.Cells.Item(29, 12).Value = "NA"
$ws.Cells.Item(29, 13).Value = "NA"
$ws.Cells.Item(30, 4).Value = 0
$ws.Cells.Item(30, 5).Value = 0
$ws.Cells.Item(30, 6).Value = 0
$ws.Cells.Item(30, 7).Value = 0
$ws.Cells.Item(30, 8).Value = 0
$ws.Cells.Item(30, 9).Value = 0
$ws.Cells.Item(30, 10).Value = 0
$ws.Cells.Item(30, 11).Value = 0
$ws.Cells.Item(30, 12).Value = 0
$ws.Cells.Item(30, 13).Value = 0
$ws.Cells.Item(31, 4).Value = 0
$ws.Cells.Item(31, 5).Value = 0
$ws.Cells.Item(31, 6).Value = 0
$ws.Cells.Item(31, 7).Value = 0
$ws.Cells.Item(31, 8).Value = 0
$ws.Cells.Item(31, 9).Value = 0
$ws.Cells.Item(31, 10).Value = 0
$ws.Cells.Item(31, 11).Value = 0
$ws.Cells.Item(31, 12).Value = 0
$ws.Cells.Item(31, 13).Value = 0
$ws.Cells.Item(32, 4).Value = -6000
$ws.Cells.Item(32, 5).Value = -4000
$ws.Cells.Item(32, 6).Value = -9000
$ws.Cells.Item(32, 7).Value = 4000
$ws.Cells.Item(32, 8).Value = -11000
$ws.Cells.Item(32, 9).Value = -11000
$ws.Cells.Item(32, 10).Value = -9000
$ws.Cells.Item(32, 11).Value = -6000
$ws.Cells.Item(32, 12).Value = -8000
$ws.Cells.Item(32, 13).Value = -7000
$ws.Cells.Item(33, 4).Value = 23000
$ws.Cells.Item(33, 5).Value = 217000
$ws.Cells.Item(33, 6).Value = 193000
$ws.Cells.Item(33, 7).Value = 4000
$ws.Cells.Item(33, 8).Value = 315000
$ws.Cells.Item(33, 9).Value = 259000
$ws.Cells.Item(33, 10).Value = 293000
$ws.Cells.Item(33, 11).Value = 93000
$ws.Cells.Item(33, 12).Value = 114000
$ws.Cells.Item(33, 13).Value = 256000
$ws.Cells.Item(34, 4).Value = 0
$ws.Cells.Item(34, 5).Value = 0
$ws.Cells.Item(34, 6).Value = 0
$ws.Cells.Item(34, 7).Value = 0
$ws.Cells.Item(34, 8).Value = 0
$ws.Cells.Item(34, 9).Value = 0
$ws.Cells.Item(34, 10).Value = 0
$ws.Cells.Item(34, 11).Value = 0
$ws.Cells.Item(34, 12).Value = 0
$ws.Cells.Item(34, 13).Value = 0
$ws.Cells.Item(35, 4).Value = 23000
$ws.Cells.Item(35, 5).Value = 217000
$ws.Cells.Item(35, 6).Value = 193000
$ws.Cells.Item(35, 7).Value = 4000
$ws.Cells.Item(35, 8).Value = 315000
$ws.Cells.Item(35, 9).Value = 259000
$ws.Cells.Item(35, 10).Value = 293000
$ws.Cells.Item(35, 11).Value = 93000
$ws.Cells.Item(35, 12).Value = 114000
$ws.Cells.Item(35, 13).Value = 256000
$ws.Cells.Item(38, 4).Value = 43465
$ws.Cells.Item(38, 5).Value = 43373
$ws.Cells.Item(38, 6).Value = 43281
$ws.Cells.Item(38, 7).Value = 43190
$ws.Cells.Item(38, 8).Value = 43100
$ws.Cells.Item(38, 9).Value = 43008
$ws.Cells.Item(38, 10).Value = 42916
$ws.Cells.Item(38, 11).Value = 42825
$ws.Cells.Item(38, 12).Value = 42735
$ws.Cells.Item(38, 13).Value = 42643
$ws.Cells.Item(41, 4).Value = 105000
$ws.Cells.Item(41, 5).Value = 174000
$ws.Cells.Item(41, 6).Value = 102000
$ws.Cells.Item(41, 7).Value = 177000
$ws.Cells.Item(41, 8).Value = 194000
$ws.Cells.Item(41, 9).Value = 144000
$ws.Cells.Item(41, 10).Value = 198000
$ws.Cells.Item(41, 11).Value = 183000
$ws.Cells.Item(41, 12).Value = 328000
$ws.Cells.Item(41, 13).Value = 1818000
$ws.Cells.Item(42, 4).Value = 1131000
$ws.Cells.Item(42, 5).Value = 1223000
$ws.Cells.Item(42, 6).Value = 1466000
$ws.Cells.Item(42, 7).Value = 1351000
$ws.Cells.Item(42, 8).Value = 1427000
$ws.Cells.Item(42, 9).Value = 1596000
$ws.Cells.Item(42, 10).Value = 1724000
$ws.Cells.Item(42, 11).Value = 1527000
$ws.Cells.Item(42, 12).Value = 1252000
$ws.Cells.Item(42, 13).Value = 1408000
$ws.Cells.Item(43, 4).Value = 366000
$ws.Cells.Item(43, 5).Value = 422000
$ws.Cells.Item(43, 6).Value = 411000
$ws.Cells.Item(43, 7).Value = 350000
$ws.Cells.Item(43, 8).Value = 682000
$ws.Cells.Item(43, 9).Value = 301000
$ws.Cells.Item(43, 10).Value = 326000
$ws.Cells.Item(43, 11).Value = 321000
$ws.Cells.Item(43, 12).Value = 302000
$ws.Cells.Item(43, 13).Value = 232000
$ws.Cells.Item(44, 4).Value = 60000
$ws.Cells.Item(44, 5).Value = 57000
$ws.Cells.Item(44, 6).Value = 57000
$ws.Cells.Item(44, 7).Value = 62000
$ws.Cells.Item(44, 8).Value = 57000
$ws.Cells.Item(44, 9).Value = 57000
$ws.Cells.Item(44, 10).Value = 52000
$ws.Cells.Item(44, 11).Value = 50000
$ws.Cells.Item(44, 12).Value = 47000
$ws.Cells.Item(44, 13).Value = 44000
$ws.Cells.Item(45, 4).Value = 125000
$ws.Cells.Item(45, 5).Value = 180000
$ws.Cells.Item(45, 6).Value = 180000
$ws.Cells.Item(45, 7).Value = 174000
$ws.Cells.Item(45, 8).Value = 133000
$ws.Cells.Item(45, 9).Value = 116000
$ws.Cells.Item(45, 10).Value = 125000
$ws.Cells.Item(45, 11).Value = 132000
$ws.Cells.Item(45, 12).Value = 121000
$ws.Cells.Item(45, 13).Value = 98000
$ws.Cells.Item(46, 4).Value = 1787000
$ws.Cells.Item(46, 5).Value = 2056000
$ws.Cells.Item(46, 6).Value = 2216000
$ws.Cells.Item(46, 7).Value = 2114000
$ws.Cells.Item(46, 8).Value = 2152000
$ws.Cells.Item(46, 9).Value = 2214000
$ws.Cells.Item(46, 10).Value = 2425000
$ws.Cells.Item(46, 11).Value = 2213000
$ws.Cells.Item(46, 12).Value = 2050000
$ws.Cells.Item(46, 13).Value = 3600000
$ws.Cells.Item(47, 4).Value = "NA"
$ws.Cells.Item(47, 5).Value = "NA"
$ws.Cells.Item(47, 6).Value = "NA"
$ws.Cells.Item(47, 7).Value = "NA"
$ws.Cells.Item(47, 8).Value = 3000
$ws.Cells.Item(47, 9).Value = "NA"
$ws.Cells.Item(47, 10).Value = "NA"
$ws.Cells.Item(47, 11).Value = "NA"
$ws.Cells.Item(47, 12).Value = "NA"
$ws.Cells.Item(47, 13).Value = "NA"
$ws.Cells.Item(48, 4).Value = 6781000
$ws.Cells.Item(48, 5).Value = 6495000
$ws.Cells.Item(48, 6).Value = 6493000
$ws.Cells.Item(48, 7).Value = 6403000
$ws.Cells.Item(48, 8).Value = 6284000
$ws.Cells.Item(48, 9).Value = 6230000
$ws.Cells.Item(48, 10).Value = 6002000
$ws.Cells.Item(48, 11).Value = 5809000
$ws.Cells.Item(48, 12).Value = 5666000
$ws.Cells.Item(48, 13).Value = 5031000
$ws.Cells.Item(49, 4).Value = 2070000
$ws.Cells.Item(49, 5).Value = 2071000
$ws.Cells.Item(49, 6).Value = 2073000
$ws.Cells.Item(49, 7).Value = 2075000
$ws.Cells.Item(49, 8).Value = 2076000
$ws.Cells.Item(49, 9).Value = 2069000
$ws.Cells.Item(49, 10).Value = 2077000
$ws.Cells.Item(49, 11).Value = 2081000
$ws.Cells.Item(49, 12).Value = 2077000
$ws.Cells.Item(49, 13).Value = "NA"
$ws.Cells.Item(50, 4).Value = 0
$ws.Cells.Item(50, 5).Value = 0
$ws.Cells.Item(50, 6).Value = 0
$ws.Cells.Item(50, 7).Value = 0
$ws.Cells.Item(50, 8).Value = 0
$ws.Cells.Item(50, 9).Value = 0
$ws.Cells.Item(50, 10).Value = 0
$ws.Cells.Item(50, 11).Value = 0
$ws.Cells.Item(50, 12).Value = 0
$ws.Cells.Item(50, 13).Value = 0
$ws.Cells.Item(51, 4).Value = 0
$ws.Cells.Item(51, 5).Value = 0
$ws.Cells.Item(51, 6).Value = 0
$ws.Cells.Item(51, 7).Value = 0
$ws.Cells.Item(51, 8).Value = 0
$ws.Cells.Item(51, 9).Value = 0
$ws.Cells.Item(51, 10).Value = 0
$ws.Cells.Item(51, 11).Value = 0
$ws.Cells.Item(51, 12).Value = 0
$ws.Cells.Item(51, 13).Value = 0
$ws.Cells.Item(52, 4).Value = 274000
$ws.Cells.Item(52, 5).Value = 271000
$ws.Cells.Item(52, 6).Value = 273000
$ws.Cells.Item(52, 7).Value = 256000
$ws.Cells.Item(52, 8).Value = 465000
$ws.Cells.Item(52, 9).Value = 226000
$ws.Cells.Item(52, 10).Value = 216000
$ws.Cells.Item(52, 11).Value = 199000
$ws.Cells.Item(52, 12).Value = 169000
$ws.Cells.Item(52, 13).Value = 68000
$ws.Cells.Item(53, 4).Value = 0
$ws.Cells.Item(53, 5).Value = 0
$ws.Cells.Item(53, 6).Value = 0
$ws.Cells.Item(53, 7).Value = 0
$ws.Cells.Item(53, 8).Value = 0
$ws.Cells.Item(53, 9).Value = 0
$ws.Cells.Item(53, 10).Value = 0
$ws.Cells.Item(53, 11).Value = 0
$ws.Cells.Item(53, 12).Value = 0
$ws.Cells.Item(53, 13).Value = 0
$ws.Cells.Item(54, 4).Value = 10912000
$ws.Cells.Item(54, 5).Value = 10893000
$ws.Cells.Item(54, 6).Value = 11055000
$ws.Cells.Item(54, 7).Value = 10848000
$ws.Cells.Item(54, 8).Value = 10746000
$ws.Cells.Item(54, 9).Value = 10739000
$ws.Cells.Item(54, 10).Value = 10720000
$ws.Cells.Item(54, 11).Value = 10302000
$ws.Cells.Item(54, 12).Value = 9962000
$ws.Cells.Item(54, 13).Value = 8699000
$ws.Cells.Item(57, 4).Value = 132000
$ws.Cells.Item(57, 5).Value = 114000
$ws.Cells.Item(57, 6).Value = 115000
$ws.Cells.Item(57, 7).Value = 102000
$ws.Cells.Item(57, 8).Value = 120000
$ws.Cells.Item(57, 9).Value = 97000
$ws.Cells.Item(57, 10).Value = 95000
$ws.Cells.Item(57, 11).Value = 95000
$ws.Cells.Item(57, 12).Value = 92000
$ws.Cells.Item(57, 13).Value = 71000
$ws.Cells.Item(58, 4).Value = 486000
$ws.Cells.Item(58, 5).Value = 345000
$ws.Cells.Item(58, 6).Value = 314000
$ws.Cells.Item(58, 7).Value = 387000
$ws.Cells.Item(58, 8).Value = 307000
$ws.Cells.Item(58, 9).Value = 334000
$ws.Cells.Item(58, 10).Value = 337000
$ws.Cells.Item(58, 11).Value = 332000
$ws.Cells.Item(58, 12).Value = 319000
$ws.Cells.Item(58, 13).Value = 275000
$ws.Cells.Item(59, 4).Value = 2324000
$ws.Cells.Item(59, 5).Value = 2429000
$ws.Cells.Item(59, 6).Value = 2607000
$ws.Cells.Item(59, 7).Value = 2503000
$ws.Cells.Item(59, 8).Value = 2259000
$ws.Cells.Item(59, 9).Value = 2334000
$ws.Cells.Item(59, 10).Value = 2542000
$ws.Cells.Item(59, 11).Value = 2408000
$ws.Cells.Item(59, 12).Value = 2124000
$ws.Cells.Item(59, 13).Value = 1777000
$ws.Cells.Item(60, 4).Value = 2942000
$ws.Cells.Item(60, 5).Value = 2888000
$ws.Cells.Item(60, 6).Value = 3036000
$ws.Cells.Item(60, 7).Value = 2992000
$ws.Cells.Item(60, 8).Value = 2686000
$ws.Cells.Item(60, 9).Value = 2765000
$ws.Cells.Item(60, 10).Value = 2974000
$ws.Cells.Item(60, 11).Value = 2835000
$ws.Cells.Item(60, 12).Value = 2535000
$ws.Cells.Item(60, 13).Value = 2123000
$ws.Cells.Item(61, 4).Value = 1617000
$ws.Cells.Item(61, 5).Value = 1684000
$ws.Cells.Item(61, 6).Value = 1998000
$ws.Cells.Item(61, 7).Value = 2062000
$ws.Cells.Item(61, 8).Value = 2262000
$ws.Cells.Item(61, 9).Value = 2367000
$ws.Cells.Item(61, 10).Value = 2469000
$ws.Cells.Item(61, 11).Value = 2531000
$ws.Cells.Item(61, 12).Value = 2645000
$ws.Cells.Item(61, 13).Value = 1861000
$ws.Cells.Item(62, 4).Value = 2602000
$ws.Cells.Item(62, 5).Value = 2530000
$ws.Cells.Item(62, 6).Value = 2428000
$ws.Cells.Item(62, 7).Value = 2355000
$ws.Cells.Item(62, 8).Value = 2338000
$ws.Cells.Item(62, 9).Value = 2116000
$ws.Cells.Item(62, 10).Value = 2013000
$ws.Cells.Item(62, 11).Value = 1922000
$ws.Cells.Item(62, 12).Value = 1851000
$ws.Cells.Item(62, 13).Value = 1851000
$ws.Cells.Item(63, 4).Value = 0
$ws.Cells.Item(63, 5).Value = 0
$ws.Cells.Item(63, 6).Value = 0
$ws.Cells.Item(63, 7).Value = 0
$ws.Cells.Item(63, 8).Value = 0
$ws.Cells.Item(63, 9).Value = 0
$ws.Cells.Item(63, 10).Value = 0
$ws.Cells.Item(63, 11).Value = 0
$ws.Cells.Item(63, 12).Value = 0
$ws.Cells.Item(63, 13).Value = 0
$ws.Cells.Item(64, 4).Value = 0
$ws.Cells.Item(64, 5).Value = 0
$ws.Cells.Item(64, 6).Value = 0
$ws.Cells.Item(64, 7).Value = 0
$ws.Cells.Item(64, 8).Value = 0
$ws.Cells.Item(64, 9).Value = 0
$ws.Cells.Item(64, 10).Value = 0
$ws.Cells.Item(64, 11).Value = 0
$ws.Cells.Item(64, 12).Value = 0
$ws.Cells.Item(64, 13).Value = 0
$ws.Cells.Item(65, 4).Value = 0
$ws.Cells.Item(65, 5).Value = 0
$ws.Cells.Item(65, 6).Value = 0
$ws.Cells.Item(65, 7).Value = 0
$ws.Cells.Item(65, 8).Value = 0
$ws.Cells.Item(65, 9).Value = 0
$ws.Cells.Item(65, 10).Value = 0
$ws.Cells.Item(65, 11).Value = 0
$ws.Cells.Item(65, 12).Value = 0
$ws.Cells.Item(65, 13).Value = 0
$ws.Cells.Item(66, 4).Value = 7161000
$ws.Cells.Item(66, 5).Value = 7102000
$ws.Cells.Item(66, 6).Value = 7462000
$ws.Cells.Item(66, 7).Value = 7409000
$ws.Cells.Item(66, 8).Value = 7286000
$ws.Cells.Item(66, 9).Value = 7248000
$ws.Cells.Item(66, 10).Value = 7456000
$ws.Cells.Item(66, 11).Value = 7288000
$ws.Cells.Item(66, 12).Value = 7031000
$ws.Cells.Item(66, 13).Value = 5835000
$ws.Cells.Item(68, 4).Value = 0
$ws.Cells.Item(68, 5).Value = 0
$ws.Cells.Item(68, 6).Value = 0
$ws.Cells.Item(68, 7).Value = 0
$ws.Cells.Item(68, 8).Value = 0
$ws.Cells.Item(68, 9).Value = 0
$ws.Cells.Item(68, 10).Value = 0
$ws.Cells.Item(68, 11).Value = 0
$ws.Cells.Item(68, 12).Value = 0
$ws.Cells.Item(68, 13).Value = 0
$ws.Cells.Item(69, 4).Value = 0
$ws.Cells.Item(69, 5).Value = 0
$ws.Cells.Item(69, 6).Value = 0
$ws.Cells.Item(69, 7).Value = 0
$ws.Cells.Item(69, 8).Value = 0
$ws.Cells.Item(69, 9).Value = 0
$ws.Cells.Item(69, 10).Value = 0
$ws.Cells.Item(69, 11).Value = 0
$ws.Cells.Item(69, 12).Value = 0
$ws.Cells.Item(69, 13).Value = 0
$ws.Cells.Item(70, 4).Value = 0
$ws.Cells.Item(70, 5).Value = 0
$ws.Cells.Item(70, 6).Value = 0
$ws.Cells.Item(70, 7).Value = 0
$ws.Cells.Item(70, 8).Value = 0
$ws.Cells.Item(70, 9).Value = 0
$ws.Cells.Item(70, 10).Value = 0
$ws.Cells.Item(70, 11).Value = 0
$ws.Cells.Item(70, 12).Value = 0
$ws.Cells.Item(70, 13).Value = 0
$ws.Cells.Item(71, 4).Value = 0
$ws.Cells.Item(71, 5).Value = 0
$ws.Cells.Item(71, 6).Value = 0
$ws.Cells.Item(71, 7).Value = 0
$ws.Cells.Item(71, 8).Value = 0
$ws.Cells.Item(71, 9).Value = 0
$ws.Cells.Item(71, 10).Value = 0
$ws.Cells.Item(71, 11).Value = 0
$ws.Cells.Item(71, 12).Value = 0
$ws.Cells.Item(71, 13).Value = 0
$ws.Cells.Item(72, 4).Value = 4534000
$ws.Cells.Item(72, 5).Value = 4550000
$ws.Cells.Item(72, 6).Value = 4373000
$ws.Cells.Item(72, 7).Value = 4219000
$ws.Cells.Item(72, 8).Value = 4193000
$ws.Cells.Item(72, 9).Value = 4406000
$ws.Cells.Item(72, 10).Value = 4182000
$ws.Cells.Item(72, 11).Value = 3629000
$ws.Cells.Item(72, 12).Value = 3568000
$ws.Cells.Item(72, 13).Value = 3488000
$ws.Cells.Item(73, 4).Value = 0
$ws.Cells.Item(73, 5).Value = 0
$ws.Cells.Item(73, 6).Value = 0
$ws.Cells.Item(73, 7).Value = 0
$ws.Cells.Item(73, 8).Value = 0
$ws.Cells.Item(73, 9).Value = 0
$ws.Cells.Item(73, 10).Value = 0
$ws.Cells.Item(73, 11).Value = 0
$ws.Cells.Item(73, 12).Value = 0
$ws.Cells.Item(73, 13).Value = 0
$ws.Cells.Item(74, 4).Value = 0
$ws.Cells.Item(74, 5).Value = 0
$ws.Cells.Item(74, 6).Value = 0
$ws.Cells.Item(74, 7).Value = 0
$ws.Cells.Item(74, 8).Value = 0
$ws.Cells.Item(74, 9).Value = 0
$ws.Cells.Item(74, 10).Value = 0
$ws.Cells.Item(74, 11).Value = 0
$ws.Cells.Item(74, 12).Value = 0
$ws.Cells.Item(74, 13).Value = 0
$ws.Cells.Item(75, 4).Value = 0
$ws.Cells.Item(75, 5).Value = 0
$ws.Cells.Item(75, 6).Value = 0
$ws.Cells.Item(75, 7).Value = 0
$ws.Cells.Item(75, 8).Value = 0
$ws.Cells.Item(75, 9).Value = 0
$ws.Cells.Item(75, 10).Value = 0
$ws.Cells.Item(75, 11).Value = 0
$ws.Cells.Item(75, 12).Value = 0
$ws.Cells.Item(75, 13).Value = 0
$ws.Cells.Item(76, 4).Value = 3751000
$ws.Cells.Item(76, 5).Value = 3791000
$ws.Cells.Item(76, 6).Value = 3593000
$ws.Cells.Item(76, 7).Value = 3439000
$ws.Cells.Item(76, 8).Value = 3460000
$ws.Cells.Item(76, 9).Value = 3491000
$ws.Cells.Item(76, 10).Value = 3264000
$ws.Cells.Item(76, 11).Value = 3014000
$ws.Cells.Item(76, 12).Value = 2931000
$ws.Cells.Item(76, 13).Value = 2864000
$ws.Cells.Item(77, 4).Value = 0
$ws.Cells.Item(77, 5).Value = 0
$ws.Cells.Item(77, 6).Value = 0
$ws.Cells.Item(77, 7).Value = 0
$ws.Cells.Item(77, 8).Value = 0
$ws.Cells.Item(77, 9).Value = 0
$ws.Cells.Item(77, 10).Value = 0
$ws.Cells.Item(77, 11).Value = 0
$ws.Cells.Item(77, 12).Value = 0
$ws.Cells.Item(77, 13).Value = 0
$ws.Cells.Item(80, 4).Value = 43465
$ws.Cells.Item(80, 5).Value = 43373
$ws.Cells.Item(80, 6).Value = 43281
$ws.Cells.Item(80, 7).Value = 43190
$ws.Cells.Item(80, 8).Value = 43100
$ws.Cells.Item(80, 9).Value = 43008
$ws.Cells.Item(80, 10).Value = 42916
$ws.Cells.Item(80, 11).Value = 42825
$ws.Cells.Item(80, 12).Value = 42735
$ws.Cells.Item(80, 13).Value = 42643
$ws.Cells.Item(81, 4).Value = 23000
$ws.Cells.Item(81, 5).Value = 217000
$ws.Cells.Item(81, 6).Value = 193000
$ws.Cells.Item(81, 7).Value = 4000
$ws.Cells.Item(81, 8).Value = 315000
$ws.Cells.Item(81, 9).Value = 259000
$ws.Cells.Item(81, 10).Value = 293000
$ws.Cells.Item(81, 11).Value = 93000
$ws.Cells.Item(81, 12).Value = 114000
$ws.Cells.Item(81, 13).Value = 256000
$ws.Cells.Item(83, 4).Value = 108000
$ws.Cells.Item(83, 5).Value = 99000
$ws.Cells.Item(83, 6).Value = 97000
$ws.Cells.Item(83, 7).Value = 94000
$ws.Cells.Item(83, 8).Value = 97000
$ws.Cells.Item(83, 9).Value = 95000
$ws.Cells.Item(83, 10).Value = 90000
$ws.Cells.Item(83, 11).Value = 90000
$ws.Cells.Item(83, 12).Value = 82000
$ws.Cells.Item(83, 13).Value = 101000
$ws.Cells.Item(84, 4).Value = 0
$ws.Cells.Item(84, 5).Value = 0
$ws.Cells.Item(84, 6).Value = 0
$ws.Cells.Item(84, 7).Value = 0
$ws.Cells.Item(84, 8).Value = 0
$ws.Cells.Item(84, 9).Value = 0
$ws.Cells.Item(84, 10).Value = 0
$ws.Cells.Item(84, 11).Value = 0
$ws.Cells.Item(84, 12).Value = 0
$ws.Cells.Item(84, 13).Value = 0
$ws.Cells.Item(85, 4).Value = 0
$ws.Cells.Item(85, 5).Value = 0
$ws.Cells.Item(85, 6).Value = 0
$ws.Cells.Item(85, 7).Value = 0
$ws.Cells.Item(85, 8).Value = 0
$ws.Cells.Item(85, 9).Value = 0
$ws.Cells.Item(85, 10).Value = 0
$ws.Cells.Item(85, 11).Value = 0
$ws.Cells.Item(85, 12).Value = 0
$ws.Cells.Item(85, 13).Value = 0
$ws.Cells.Item(86, 4).Value = 0
$ws.Cells.Item(86, 5).Value = 0
$ws.Cells.Item(86, 6).Value = 0
$ws.Cells.Item(86, 7).Value = 0
$ws.Cells.Item(86, 8).Value = 0
$ws.Cells.Item(86, 9).Value = 0
$ws.Cells.Item(86, 10).Value = 0
$ws.Cells.Item(86, 11).Value = 0
$ws.Cells.Item(86, 12).Value = 0
$ws.Cells.Item(86, 13).Value = 0
$ws.Cells.Item(87, 4).Value = 0
$ws.Cells.Item(87, 5).Value = 0
$ws.Cells.Item(87, 6).Value = 0
$ws.Cells.Item(87, 7).Value = 0
$ws.Cells.Item(87, 8).Value = 0
$ws.Cells.Item(87, 9).Value = 0
$ws.Cells.Item(87, 10).Value = 0
$ws.Cells.Item(87, 11).Value = 0
$ws.Cells.Item(87, 12).Value = 0
$ws.Cells.Item(87, 13).Value = 0
$ws.Cells.Item(88, 4).Value = 0
$ws.Cells.Item(88, 5).Value = 0
$ws.Cells.Item(88, 6).Value = 0
$ws.Cells.Item(88, 7).Value = 0
$ws.Cells.Item(88, 8).Value = 0
$ws.Cells.Item(88, 9).Value = 0
$ws.Cells.Item(88, 10).Value = 0
$ws.Cells.Item(88, 11).Value = 0
$ws.Cells.Item(88, 12).Value = 0
$ws.Cells.Item(88, 13).Value = 0
$ws.Cells.Item(89, 4).Value = 209000
$ws.Cells.Item(89, 5).Value = 260000
$ws.Cells.Item(89, 6).Value = 420000
$ws.Cells.Item(89, 7).Value = 306000
$ws.Cells.Item(89, 8).Value = 233000
$ws.Cells.Item(89, 9).Value = 273000
$ws.Cells.Item(89, 10).Value = 614000
$ws.Cells.Item(89, 11).Value = 470000
$ws.Cells.Item(89, 12).Value = 180000
$ws.Cells.Item(89, 13).Value = 307000
$ws.Cells.Item(91, 4).Value = -29000
$ws.Cells.Item(91, 5).Value = -20000
$ws.Cells.Item(91, 6).Value = -27000
$ws.Cells.Item(91, 7).Value = -29000
$ws.Cells.Item(91, 8).Value = -26000
$ws.Cells.Item(91, 9).Value = -25000
$ws.Cells.Item(91, 10).Value = -25000
$ws.Cells.Item(91, 11).Value = -432000
$ws.Cells.Item(91, 12).Value = -169000
$ws.Cells.Item(91, 13).Value = -169000
$ws.Cells.Item(92, 4).Value = 0
$ws.Cells.Item(92, 5).Value = 0
$ws.Cells.Item(92, 6).Value = 0
$ws.Cells.Item(92, 7).Value = 0
$ws.Cells.Item(92, 8).Value = 0
$ws.Cells.Item(92, 9).Value = 0
$ws.Cells.Item(92, 10).Value = 0
$ws.Cells.Item(92, 11).Value = 0
$ws.Cells.Item(92, 12).Value = 0
$ws.Cells.Item(92, 13).Value = 0
$ws.Cells.Item(93, 4).Value = 0
$ws.Cells.Item(93, 5).Value = 0
$ws.Cells.Item(93, 6).Value = 0
$ws.Cells.Item(93, 7).Value = 0
$ws.Cells.Item(93, 8).Value = 0
$ws.Cells.Item(93, 9).Value = 0
$ws.Cells.Item(93, 10).Value = 0
$ws.Cells.Item(93, 11).Value = 0
$ws.Cells.Item(93, 12).Value = 0
$ws.Cells.Item(93, 13).Value = 0
$ws.Cells.Item(94, 4).Value = -298000
$ws.Cells.Item(94, 5).Value = 135000
$ws.Cells.Item(94, 6).Value = -304000
$ws.Cells.Item(94, 7).Value = -164000
$ws.Cells.Item(94, 8).Value = 13000
$ws.Cells.Item(94, 9).Value = -171000
$ws.Cells.Item(94, 10).Value = -483000
$ws.Cells.Item(94, 11).Value = -488000
$ws.Cells.Item(94, 12).Value = -1981000
$ws.Cells.Item(94, 13).Value = -51000
$ws.Cells.Item(96, 4).Value = -40000
$ws.Cells.Item(96, 5).Value = -39000
$ws.Cells.Item(96, 6).Value = -40000
$ws.Cells.Item(96, 7).Value = -39000
$ws.Cells.Item(96, 8).Value = -37000
$ws.Cells.Item(96, 9).Value = -37000
$ws.Cells.Item(96, 10).Value = -37000
$ws.Cells.Item(96, 11).Value = -37000
$ws.Cells.Item(96, 12).Value = -34000
$ws.Cells.Item(96, 13).Value = -34000
$ws.Cells.Item(97, 4).Value = 0
$ws.Cells.Item(97, 5).Value = 0
$ws.Cells.Item(97, 6).Value = 0
$ws.Cells.Item(97, 7).Value = 0
$ws.Cells.Item(97, 8).Value = 0
$ws.Cells.Item(97, 9).Value = 0
$ws.Cells.Item(97, 10).Value = 0
$ws.Cells.Item(97, 11).Value = 0
$ws.Cells.Item(97, 12).Value = 0
$ws.Cells.Item(97, 13).Value = 0
$ws.Cells.Item(98, 4).Value = 0
$ws.Cells.Item(98, 5).Value = 0
$ws.Cells.Item(98, 6).Value = 0
$ws.Cells.Item(98, 7).Value = 0
$ws.Cells.Item(98, 8).Value = 0
$ws.Cells.Item(98, 9).Value = 0
$ws.Cells.Item(98, 10).Value = 0
$ws.Cells.Item(98, 11).Value = 0
$ws.Cells.Item(98, 12).Value = 0
$ws.Cells.Item(98, 13).Value = 0
$ws.Cells.Item(99, 4).Value = 0
$ws.Cells.Item(99, 5).Value = 0
$ws.Cells.Item(99, 6).Value = 0
$ws.Cells.Item(99, 7).Value = 0
$ws.Cells.Item(99, 8).Value = 0
$ws.Cells.Item(99, 9).Value = 0
$ws.Cells.Item(99, 10).Value = 0
$ws.Cells.Item(99, 11).Value = 0
$ws.Cells.Item(99, 12).Value = 0
$ws.Cells.Item(99, 13).Value = 0
$ws.Cells.Item(100, 4).Value = 19000
$ws.Cells.Item(100, 5).Value = -319000
$ws.Cells.Item(100, 6).Value = -191000
$ws.Cells.Item(100, 7).Value = -156000
$ws.Cells.Item(100, 8).Value = -193000
$ws.Cells.Item(100, 9).Value = -156000
$ws.Cells.Item(100, 10).Value = -116000
$ws.Cells.Item(100, 11).Value = -127000
$ws.Cells.Item(100, 12).Value = 311000
$ws.Cells.Item(100, 13).Value = 1481000
$ws.Cells.Item(101, 4).Value = 0
$ws.Cells.Item(101, 5).Value = 0
$ws.Cells.Item(101, 6).Value = 0
$ws.Cells.Item(101, 7).Value = 0
$ws.Cells.Item(101, 8).Value = 0
$ws.Cells.Item(101, 9).Value = 0
$ws.Cells.Item(101, 10).Value = 0
$ws.Cells.Item(101, 11).Value = 0
$ws.Cells.Item(101, 12).Value = 0
$ws.Cells.Item(101, 13).Value = 0
$ws.Cells.Item(102, 4).Value = -70000
$ws.Cells.Item(102, 5).Value = 76000
$ws.Cells.Item(102, 6).Value = -75000
$ws.Cells.Item(102, 7).Value = -14000
$ws.Cells.Item(102, 8).Value = 53000
$ws.Cells.Item(102, 9).Value = -54000
$ws.Cells.Item(102, 10).Value = 15000
$ws.Cells.Item(102, 11).Value = -145000
$ws.Cells.Item(102, 12).Value = -1490000
$ws.Cells.Item(102, 13).Value = 1737000
